# reverse engineered checkPassword on hw3_dry
#
# Fills in the blanks of the checkPassword() / checkPasswordAux() code
# listing with the reverse-engineered answers.
#
# Word's Font.Color is an OLE_COLOR int stored 0x00BBGGRR (i.e. the byte
# order is reversed from the RRGGBB hex used in the document's <w:color>
# attributes), hence the helper below.

$d = $word.ActiveDocument

function HexToWordColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$colorDefault = HexToWordColor "D4D4D4"   # plain code text (VS Code dark+ fg)
$colorString  = HexToWordColor "CE9178"   # string-literal orange

# Monotonically increasing search cursor: placeholders like "____" repeat
# throughout the listing, so every lookup is scoped to [cursor, end) and
# resolves to the next occurrence in document order.
$cursor = 0

function FindNext([string]$text) {
    $rng = $d.Range($cursor, $d.Content.End)
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    if (-not $rng.Find.Found) {
        throw "Could not find text: $text"
    }
    $script:cursor = $rng.End
    return $rng
}

# Re-applies a color to a just-inserted sub-range so the save pass keeps
# it as its own <w:r> instead of silently coalescing it with a
# same-formatted neighbour run.
function SplitColor($range, $color) {
    $range.Font.Color = 255
    $range.Font.Color = $color
}

# --- 1) while(_______){  ->  while(s != NULL){ -------------------------
$r = FindNext "_______"
$r.Text = "s != NULL"

# --- 2) unsigned long x = *s - ______;  ->  unsigned long x = *s – 'a'; -
$r = FindNext " x = *s -"
$r.Text = " x = *s "
$r.InsertAfter([string][char]0x2013)   # "–" en dash
$cursor = $r.End
SplitColor ($d.Range($r.End - 1, $r.End)) $colorDefault

$r = FindNext " ______"
$r.Text = " "
$r.InsertAfter("'a'")
$cursor = $r.End
SplitColor ($d.Range($r.End - 3, $r.End)) $colorString

# --- 3) if(x>____){  ->  if(x > 25){ ------------------------------------
$r = FindNext "(x>"
$r.Text = "(x"
$r.InsertAfter(" ")
$r.InsertAfter(">")
$r.InsertAfter(" ")
$cursor = $r.End
# back-fill formatting from the tail forward so earlier offsets stay valid
SplitColor ($d.Range($r.End - 1, $r.End)) $colorDefault       # " "
SplitColor ($d.Range($r.End - 2, $r.End - 1)) $colorDefault   # ">"
SplitColor ($d.Range($r.End - 3, $r.End - 2)) $colorDefault   # " "

$r = FindNext "____"
$r.Text = "25"

# --- 4) if(y > _____){  ->  if(y > ~(x)){ -------------------------------
$r = FindNext "(y > "
$r.Text = "(y >"

$r = FindNext "_____"
$r.Text = " "
$r.InsertAfter("~(x)")
$cursor = $r.End
SplitColor ($d.Range($r.End - 4, $r.End)) $colorDefault

# --- 5) y = __________;  ->  y = 26y + *s; ------------------------------
$r = FindNext "         y = "
$r.Text = "         y ="

$r = FindNext "__________"
$r.Text = " 26y + *s"

# --- 6) return ____________;  ->  return y == hash; ---------------------
$r = FindNext "____________"
$r.Text = "y == hash"

Write-Output "done"
